$d = $word.ActiveDocument

# The "Date" paragraph currently reads "15:54 17 September 2019" and is
# built from the separate runs: "15:54", " ", "17", " ", "September",
# " ", "2019". We need to drop the leading "15:54" and the following
# space, leaving "17 September 2019" split across its original runs
# (17 / " " / September / " " / 2019) untouched.

$datePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Style.NameLocal -eq "Date") {
        $datePara = $candidate
        break
    }
}

$paraRange = $datePara.Range
# Exclude the trailing paragraph mark from the replacement range.
$textRange = $d.Range($paraRange.Start, $paraRange.End - 1)

$flatOpc = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p>' +
    '<w:r><w:t xml:space="preserve">17</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">September</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">2019</w:t></w:r>' +
    '</w:p></w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$textRange.InsertXML($flatOpc)
